$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $savedStyle
}

Set-TextValue "D2" "60.767.39"
Set-TextValue "E2" "  +0.36%  "
Set-TextValue "D3" "2.655.64"
Set-TextValue "E3" "  +2.03%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "578.48"
Set-TextValue "E5" "  +0.94%  "
Set-TextValue "D6" "145.22"
Set-TextValue "E6" "  +1.89%  "
Set-TextValue "D7" "0.997"
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "E8" "  -0.42%  "
Set-TextValue "D9" "6.56"
Set-TextValue "E9" "  +0.10%  "
Set-TextValue "E10" "  +1.30%  "
Set-TextValue "D11" "0.375"
Set-TextValue "E11" "  +2.92%  "
Set-TextValue "E12" "  +2.50%  "
Set-TextValue "D13" "3.127.32"
Set-TextValue "E13" "  +2.03%  "
Set-TextValue "D14" "25.75"
Set-TextValue "E14" "  +10.71%  "
Set-TextValue "D15" "60.744.81"
Set-TextValue "E15" "  +0.25%  "
Set-TextValue "E16" "  +1.81%  "
Set-TextValue "D17" "2.664.06"
Set-TextValue "E17" "  +1.86%  "
Set-TextValue "D18" "11.55"
Set-TextValue "E18" "  +1.41%  "
Set-TextValue "E19" "  +1.77%  "
Set-TextValue "D20" "350.96"
Set-TextValue "E20" "  +1.15%  "
Set-TextValue "D21" "6.95"
Set-TextValue "E21" "  -0.46%  "
Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  +0.10%  "
Set-TextValue "D23" "0.533"
Set-TextValue "E23" "  +0.71%  "
Set-TextValue "D24" "64.01"
Set-TextValue "E24" "  +0.86%  "
Set-TextValue "D25" "0.997"
Set-TextValue "E25" "  +0.09%  "
Set-TextValue "E26" "  +1.92%  "
Set-TextValue "E27" "  +5.47%  "
Set-TextValue "D28" "1.99"
Set-TextValue "E28" "  +9.36%  "
Set-TextValue "E29" "  +2.99%  "
Set-TextValue "D30" "6.75"
Set-TextValue "E30" "  +5.73%  "
Set-TextValue "D31" "168.03"
Set-TextValue "E31" "  +4.75%  "
Set-TextValue "E32" "  +0.09%  "
Set-TextValue "D33" "19.91"
Set-TextValue "E33" "  +1.98%  "
Set-TextValue "D34" "1.06"
Set-TextValue "E34" "  +9.33%  "
Set-TextValue "E35" "  +5.43%  "
Set-TextValue "D36" "1.32"
Set-TextValue "E36" "  +8.38%  "
Set-TextValue "D37" "1.64"
Set-TextValue "E37" "  +2.26%  "
Set-TextValue "D38" "327.49"
Set-TextValue "E38" "  +10.84%  "
Set-TextValue "E39" "  +4.76%  "
Set-TextValue "D40" "38.41"
Set-TextValue "E40" "  +1.69%  "
Set-TextValue "D41" "0.881"
Set-TextValue "E41" "  +4.01%  "
Set-TextValue "E42" "  +6.91%  "
Set-TextValue "D43" "20.64"
Set-TextValue "E43" "  +4.43%  "
Set-TextValue "B44" "Aave"
Set-TextValue "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "134.41"
Set-TextValue "E44" "  -2.79%  "
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "20.69"
Set-TextValue "E45" "  +4.04%  "
Set-TextValue "D46" "0.0998"
Set-TextValue "E46" "  +1.39%  "
Set-TextValue "D47" "0.614"
Set-TextValue "E47" "  +0.93%  "
Set-TextValue "D48" "0.0558"
Set-TextValue "E48" "  +2.44%  "
Set-TextValue "E49" "  +0.27%  "
Set-TextValue "D50" "0.0246"
Set-TextValue "E50" "  +2.26%  "
Set-TextValue "D51" "2.139.69"
Set-TextValue "E51" "  +5.60%  "
